$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C19").Value = 0
